$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2023 December (row 5): Closed Issues 51 -> 53, Opened Issues 9 -> 7
$ws.Range("B5").Value = 53
$ws.Range("C5").Value = 7

# 2024 January (row 6): Closed Issues 52 -> 53, Opened Issues 9 -> 8
$ws.Range("B6").Value = 53
$ws.Range("C6").Value = 8

# 2024 February (row 7): Closed Issues 31 -> 33, Opened Issues 17 -> 15
$ws.Range("B7").Value = 33
$ws.Range("C7").Value = 15

# 2024 July (row 12): Closed Issues 33 -> 34, Opened Issues 20 -> 19
$ws.Range("B12").Value = 34
$ws.Range("C12").Value = 19

# 2024 August (row 13): Closed Issues 14 -> 20, Opened Issues 19 -> 22
$ws.Range("B13").Value = 20
$ws.Range("C13").Value = 22
